$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as literal text (e.g. "1.00", "10.68")
# rather than numbers, so trailing zeros and grouping survive untouched.
# Assigning a plain numeric-looking string to .Value would otherwise let
# Excel re-interpret it as a number, so for cells whose new price text
# parses as a number we briefly force a text number format, then restore
# the default "Normal" style once the literal text is committed.
$textCells = @("D5", "D6", "D10", "D12", "D13", "D14", "D17", "D18", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.013.82"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.730.57"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "616.25"
$ws.Range("E5").Value = "  +3.76%  "
$ws.Range("D6").Value = "183.73"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("D7").Value = "3.727.97"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "0.722"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("E11").Value = "  -6.78%  "
$ws.Range("D12").Value = "56.97"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "0.0000294"
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("D14").Value = "10.68"
$ws.Range("E14").Value = "  -4.60%  "
$ws.Range("D15").Value = "4.325.33"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "3.738.26"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "19.44"
$ws.Range("E17").Value = "  -3.59%  "
$ws.Range("D18").Value = "13.09"
$ws.Range("E18").Value = "  -3.70%  "
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("D21").Value = "68.891.37"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "414.13"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "4.67"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "89.86"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "3.06"
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("D26").Value = "12.83"
$ws.Range("E26").Value = "  -4.97%  "
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("D28").Value = "6.07"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "9.67"
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("D31").Value = "33.28"
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("D32").Value = "7.36"
$ws.Range("E32").Value = "  -14.91%  "
$ws.Range("D33").Value = "12.77"
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("D34").Value = "0.122"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "621.69"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "44.49"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").Value = "66.29"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "0.0₃0880"
$ws.Range("E38").Value = "  -8.87%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "0.405"
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.142"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").Value = "0.0444"
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("D45").Value = "2.64"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.140"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "9.24"
$ws.Range("E47").Value = "  -6.96%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.825.50"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  -17.62%  "
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("E51").Value = "  -7.01%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
